$d = $word.ActiveDocument

# 1) Update the "Generated:" date line.
[void]$d.Content.Find.Execute("Generated: 09 Dec 2025", $true, $false, $false, $false, $false, $true, 1, $false, "Generated: 11 Dec 2025", 2)

# 2) Rewrite the activity-log paragraph with the new chronological narrative
#    and "Current Status" bullet list.
#    That paragraph's single run stores its text as one <w:t> with <w:br/>
#    elements nested directly inside it (instead of as sibling run children),
#    so Word's normal text layer reports the paragraph as empty and
#    Find/Replace cannot locate anything inside it. We therefore locate the
#    paragraph via its raw OOXML and rebuild it wholesale through
#    Range.InsertXML, preserving the paragraph's original identity
#    attributes (paraId/textId/rsid*).
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.XML() -like "*rewritten version of the activity log*") {
        $targetPara = $candidate
        break
    }
}
if ($targetPara -eq $null) {
    $targetPara = $d.Paragraphs.Item(4)
}

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1A7F9F88" w14:textId="20BBCB5E" w:rsidR="00FD11D2" w:rsidRDefault="005D753A" w:rsidP="005D753A"><w:r><w:t>Here is the rewritten activity log in a smooth, chronological project narrative:<w:br/><w:br/>The AI document generation functionality was successfully tested and validated to ensure that its format and structure met the required specifications. Concurrently, we began integrating AI Document Generation with our GitHub repository, setting up a complete structure and branch workflow.<w:br/><w:br/>Next, we turned our attention to creating the Personal Employee Page (PEP), which involved designing both the UI and backend endpoint integration with the database. This task was completed promptly, allowing us to move forward with planning the integration of AI and polishing the backend further.<w:br/><w:br/>We then refined the API structure, improved service, controller, and data flow, and prepared for the integration of AI Document Generation. After retesting the project and task workflow to ensure stability, we fixed bugs related to data display, enhanced UX (User Experience), and validated forms. Furthermore, we successfully connected the backend, frontend, and PostgreSQL database.<w:br/><w:br/>Subsequent tasks included adding Project Form and Task Form, testing interaction flow between them, and fixing display bugs related to billable assignments, due dates, and form inputs. We also integrated these forms with the backend and ensured they received input correctly.<w:br/><w:br/>To further enhance our project, we began detailing the flow of interactions between Projects, Tasks, and Activities, as well as integrating Tailwind CSS into the frontend. Meanwhile, we improved web display, developed task and activity pages, started designing interaction flows, and initiated API structure planning using Go on the backend.<w:br/><w:br/>Additionally, we started building reporting functionality and handling file XER. To optimize our workflow, we adjusted layout, project structure, and work environment settings. Finally, we set up initial frontend and backend configurations to ensure a solid foundation for future development.<w:br/><w:br/>**Current Status:**<w:br/><w:br/>* Completed tasks:<w:br/></w:t><w:tab/><w:t>+ Tested AI document generation<w:br/></w:t><w:tab/><w:t>+ Integrated AI Document Generation with GitHub repository<w:br/></w:t><w:tab/><w:t>+ Designed Personal Employee Page UI and backend endpoint integration<w:br/></w:t><w:tab/><w:t>+ Refined API structure and service flow<w:br/></w:t><w:tab/><w:t>+ Fixed data display bugs<w:br/></w:t><w:tab/><w:t>+ Connected backend, frontend, and PostgreSQL database<w:br/></w:t><w:tab/><w:t>+ Added Project Form and Task Form<w:br/></w:t><w:tab/><w:t>+ Enhanced UX and validated forms<w:br/></w:t><w:tab/><w:t>+ Integrated Tailwind CSS into frontend<w:br/>* Pending tasks:<w:br/></w:t><w:tab/><w:t>+ Complete AI integration with backend<w:br/></w:t><w:tab/><w:t>+ Finalize API structure planning using Go on the backend<w:br/></w:t><w:tab/><w:t>+ Develop reporting functionality and file XER handling<w:br/>* Blockers: None identified at this time<w:br/>* Risks: Potential integration complexities between AI and backend functionalities</w:t></w:r></w:p>'
[void]$targetPara.Range.InsertXML($newParaXml)
